# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 347-348 (pushing the existing
# rows 347-375 down to 349-377), matching the new dataset snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 347; this shifts the old
# rows 347:375 down to 349:377 and extends the sheet dimension to R377.
$ws.Rows("347:348").Insert()

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D347:D348").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 347
$ws.Range("A347").Value = 5
$ws.Range("B347").Value = "Macroferia Regional de Talca"
$ws.Range("C347").Value = "Maule"
$ws.Range("D347").Value = 44783
$ws.Range("E347").Value = 7
$ws.Range("F347").Value = 100112023
$ws.Range("G347").Value = "Brócoli"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 3000
$ws.Range("K347").Value = 1000
$ws.Range("L347").Value = 1000
$ws.Range("M347").Value = 1000
$ws.Range("N347").Value = "$/unidad"
$ws.Range("O347").Value = "Región del Maule"
$ws.Range("P347").Value = 1000
$ws.Range("Q347").Value = 1
$ws.Range("R347").Value = "Hortaliza"

# New row 348
$ws.Range("A348").Value = 5
$ws.Range("B348").Value = "Macroferia Regional de Talca"
$ws.Range("C348").Value = "Maule"
$ws.Range("D348").Value = 44783
$ws.Range("E348").Value = 7
$ws.Range("F348").Value = 100112023
$ws.Range("G348").Value = "Brócoli"
$ws.Range("H348").Value = "Sin especificar"
$ws.Range("I348").Value = "Segunda"
$ws.Range("J348").Value = 3000
$ws.Range("K348").Value = 800
$ws.Range("L348").Value = 800
$ws.Range("M348").Value = 800
$ws.Range("N348").Value = "$/unidad"
$ws.Range("O348").Value = "Región del Maule"
$ws.Range("P348").Value = 800
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = "Hortaliza"
